$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 164
$ws.Range("J2").Value = 722
$ws.Range("K2").Value = 4
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 136
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 73
$ws.Range("T2").Value = 117
$ws.Range("V2").Value = 1135
$ws.Range("X2").Value = 1106
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 10
